# Scheduled runner update: refresh market-price derived columns
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit) across
# several crafting-job sheets with newly pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 39421.75
$ws.Range("I62").Value = 103628.9
$ws.Range("J62").Value = 3751.111
$ws.Range("K62").Value = 103628.9
$ws.Range("L62").Value = 3751.111
$ws.Range("M62").Value = -103004.9
$ws.Range("N62").Value = -4999.111
$ws.Range("H65").Value = 39421.75
$ws.Range("I65").Value = 103628.9
$ws.Range("J65").Value = 3751.111
$ws.Range("K65").Value = 518144.5
$ws.Range("L65").Value = 18755.555
$ws.Range("M65").Value = -515024.5
$ws.Range("N65").Value = -24995.555
$ws.Range("H80").Value = 458232.75
$ws.Range("I80").Value = 933.2222
$ws.Range("J80").Value = 972694.75
$ws.Range("K80").Value = 2799.6666
$ws.Range("L80").Value = 2918084.25
$ws.Range("M80").Value = -1801.6666
$ws.Range("N80").Value = -2920080.25
$ws.Range("H83").Value = 458232.75
$ws.Range("I83").Value = 933.2222
$ws.Range("J83").Value = 972694.75
$ws.Range("K83").Value = 8398.9998
$ws.Range("L83").Value = 8754252.75
$ws.Range("M83").Value = -3406.9998
$ws.Range("N83").Value = -8764236.75
$ws.Range("H92").Value = 1155.6666
$ws.Range("I92").Value = 1271.6111
$ws.Range("J92").Value = 460
$ws.Range("K92").Value = 1271.6111
$ws.Range("L92").Value = 460
$ws.Range("M92").Value = -23.61110000000008
$ws.Range("N92").Value = -2956
$ws.Range("H106").Value = 3820.1428
$ws.Range("I106").Value = 3718.2
$ws.Range("K106").Value = 3718.2
$ws.Range("M106").Value = -3087.2
$ws.Range("H112").Value = 1788.75
$ws.Range("I112").Value = 1275
$ws.Range("J112").Value = 1862.1428
$ws.Range("K112").Value = 3825
$ws.Range("L112").Value = 5586.428400000001
$ws.Range("M112").Value = -2717
$ws.Range("N112").Value = -7802.428400000001
$ws.Range("H116").Value = 3003.318
$ws.Range("I116").Value = 2767.0625
$ws.Range("J116").Value = 3633.3333
$ws.Range("K116").Value = 2767.0625
$ws.Range("L116").Value = 3633.3333
$ws.Range("M116").Value = 674.9375
$ws.Range("N116").Value = -10517.3333
$ws.Range("H130").Value = 28052
$ws.Range("J130").Value = 28052
$ws.Range("L130").Value = 28052
$ws.Range("N130").Value = -38092
$ws.Range("H132").Value = 1097.847
$ws.Range("I132").Value = 1118.5325
$ws.Range("K132").Value = 3355.5975
$ws.Range("M132").Value = -825.5974999999999
$ws.Range("H135").Value = 763.68335
$ws.Range("I135").Value = 425.64706
$ws.Range("J135").Value = 2679.2222
$ws.Range("K135").Value = 3830.82354
$ws.Range("L135").Value = 24112.9998
$ws.Range("M135").Value = -1295.82354
$ws.Range("N135").Value = -29182.9998
$ws.Range("H138").Value = 1515.89
$ws.Range("I138").Value = 804.9761999999999
$ws.Range("J138").Value = 2030.6897
$ws.Range("K138").Value = 2414.9286
$ws.Range("L138").Value = 6092.0691
$ws.Range("M138").Value = 2725.0714
$ws.Range("N138").Value = -16372.0691
$ws.Range("H141").Value = 2175.578
$ws.Range("I141").Value = 852.5263
$ws.Range("J141").Value = 9357.857
$ws.Range("K141").Value = 2557.5789
$ws.Range("L141").Value = 28073.571
$ws.Range("M141").Value = 2622.4211
$ws.Range("N141").Value = -38433.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22690.918
$ws.Range("I32").Value = 21994.521
$ws.Range("J32").Value = 25517.47
$ws.Range("K32").Value = 21994.521
$ws.Range("L32").Value = 25517.47
$ws.Range("M32").Value = -21707.521
$ws.Range("N32").Value = -26091.47
$ws.Range("H61").Value = 882.44183
$ws.Range("I61").Value = 791.1539
$ws.Range("J61").Value = 1772.5
$ws.Range("K61").Value = 791.1539
$ws.Range("L61").Value = 1772.5
$ws.Range("M61").Value = -579.1539
$ws.Range("N61").Value = -2196.5
$ws.Range("H98").Value = 23568.334
$ws.Range("J98").Value = 23568.334
$ws.Range("L98").Value = 23568.334
$ws.Range("N98").Value = -29558.334
$ws.Range("H136").Value = 882.44183
$ws.Range("I136").Value = 791.1539
$ws.Range("J136").Value = 1772.5
$ws.Range("K136").Value = 2373.4617
$ws.Range("L136").Value = 5317.5
$ws.Range("M136").Value = 176.5383000000002
$ws.Range("N136").Value = -10417.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 514.1667
$ws.Range("I80").Value = 285
$ws.Range("J80").Value = 579.6429000000001
$ws.Range("K80").Value = 285
$ws.Range("L80").Value = 579.6429000000001
$ws.Range("M80").Value = 713
$ws.Range("N80").Value = -2575.6429
$ws.Range("H83").Value = 514.1667
$ws.Range("I83").Value = 285
$ws.Range("J83").Value = 579.6429000000001
$ws.Range("K83").Value = 1425
$ws.Range("L83").Value = 2898.2145
$ws.Range("M83").Value = 3567
$ws.Range("N83").Value = -12882.2145
$ws.Range("H86").Value = 2517.3333
$ws.Range("I86").Value = 2319.476
$ws.Range("J86").Value = 3902.3333
$ws.Range("K86").Value = 2319.476
$ws.Range("L86").Value = 3902.3333
$ws.Range("M86").Value = -1196.476
$ws.Range("N86").Value = -6148.3333
$ws.Range("H89").Value = 2517.3333
$ws.Range("I89").Value = 2319.476
$ws.Range("J89").Value = 3902.3333
$ws.Range("K89").Value = 11597.38
$ws.Range("L89").Value = 19511.6665
$ws.Range("M89").Value = -5981.380000000001
$ws.Range("N89").Value = -30743.6665
$ws.Range("H107").Value = 9570.0625
$ws.Range("I107").Value = 1029.1818
$ws.Range("K107").Value = 1029.1818
$ws.Range("M107").Value = 890.8181999999999
$ws.Range("H122").Value = 36960
$ws.Range("J122").Value = 36960
$ws.Range("L122").Value = 36960
$ws.Range("N122").Value = -46760
$ws.Range("H134").Value = 16182.479
$ws.Range("I134").Value = 1458.1754
$ws.Range("J134").Value = 86122.914
$ws.Range("K134").Value = 4374.5262
$ws.Range("L134").Value = 258368.742
$ws.Range("M134").Value = -1839.5262
$ws.Range("N134").Value = -263438.742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 69257.5
$ws.Range("J23").Value = 69257.5
$ws.Range("L23").Value = 69257.5
$ws.Range("N23").Value = -69737.5
$ws.Range("H27").Value = 69257.5
$ws.Range("J27").Value = 69257.5
$ws.Range("L27").Value = 69257.5
$ws.Range("N27").Value = -69641.5
$ws.Range("H31").Value = 2095499.4
$ws.Range("I31").Value = 2370715
$ws.Range("J31").Value = 3860
$ws.Range("K31").Value = 2370715
$ws.Range("L31").Value = 3860
$ws.Range("M31").Value = -2370420
$ws.Range("N31").Value = -4450
$ws.Range("H34").Value = 2095499.4
$ws.Range("I34").Value = 2370715
$ws.Range("J34").Value = 3860
$ws.Range("K34").Value = 2370715
$ws.Range("L34").Value = 3860
$ws.Range("M34").Value = -2370513
$ws.Range("N34").Value = -4264
$ws.Range("H58").Value = 3907.3684
$ws.Range("I58").Value = 1352.75
$ws.Range("J58").Value = 8286.714
$ws.Range("K58").Value = 1352.75
$ws.Range("L58").Value = 8286.714
$ws.Range("M58").Value = -1149.75
$ws.Range("N58").Value = -8692.714
$ws.Range("H105").Value = 2799
$ws.Range("I105").Value = 2741.4285
$ws.Range("J105").Value = 2933.3333
$ws.Range("K105").Value = 2741.4285
$ws.Range("L105").Value = 2933.3333
$ws.Range("M105").Value = -994.4285
$ws.Range("N105").Value = -6427.3333
$ws.Range("H132").Value = 1483.395
$ws.Range("I132").Value = 869.6875
$ws.Range("J132").Value = 2376.0605
$ws.Range("K132").Value = 2609.0625
$ws.Range("L132").Value = 7128.181500000001
$ws.Range("M132").Value = -79.0625
$ws.Range("N132").Value = -12188.1815
$ws.Range("H134").Value = 1130.0588
$ws.Range("I134").Value = 1044.8
$ws.Range("J134").Value = 1366.8889
$ws.Range("K134").Value = 3134.4
$ws.Range("L134").Value = 4100.6667
$ws.Range("M134").Value = -599.3999999999996
$ws.Range("N134").Value = -9170.6667
$ws.Range("H136").Value = 3907.3684
$ws.Range("I136").Value = 1352.75
$ws.Range("J136").Value = 8286.714
$ws.Range("K136").Value = 4058.25
$ws.Range("L136").Value = 24860.142
$ws.Range("M136").Value = -1508.25
$ws.Range("N136").Value = -29960.142
$ws.Range("H141").Value = 58599.8
$ws.Range("J141").Value = 61461.31
$ws.Range("L141").Value = 61461.31
$ws.Range("N141").Value = -71821.31

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 400
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 450
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 1350
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -3846

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 942.625
$ws.Range("I113").Value = 315.36365
$ws.Range("J113").Value = 2322.6
$ws.Range("K113").Value = 315.36365
$ws.Range("L113").Value = 2322.6
$ws.Range("M113").Value = 1854.63635
$ws.Range("N113").Value = -6662.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4447206
$ws.Range("I7").Value = 3248.9333
$ws.Range("K7").Value = 3248.9333
$ws.Range("M7").Value = -3136.9333
$ws.Range("H16").Value = 5394.643
$ws.Range("I16").Value = 7591.647
$ws.Range("K16").Value = 7591.647
$ws.Range("M16").Value = -7421.647
$ws.Range("H126").Value = 4447206
$ws.Range("I126").Value = 3248.9333
$ws.Range("K126").Value = 9746.7999
$ws.Range("M126").Value = -7276.7999
$ws.Range("H132").Value = 2034.849
$ws.Range("I132").Value = 1888.1052
$ws.Range("J132").Value = 2406.6
$ws.Range("K132").Value = 5664.3156
$ws.Range("L132").Value = 7219.799999999999
$ws.Range("M132").Value = -3134.3156
$ws.Range("N132").Value = -12279.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 27665
$ws.Range("J110").Value = 27665
$ws.Range("L110").Value = 27665
$ws.Range("N110").Value = -35845
$ws.Range("H132").Value = 1395.7693
$ws.Range("I132").Value = 1355.7333
$ws.Range("J132").Value = 1450.3636
$ws.Range("K132").Value = 4067.199900000001
$ws.Range("L132").Value = 4351.0908
$ws.Range("M132").Value = -1537.199900000001
$ws.Range("N132").Value = -9411.0908
$ws.Range("H136").Value = 1158.5714
$ws.Range("I136").Value = 1303
$ws.Range("J136").Value = 725.2857
$ws.Range("K136").Value = 3909
$ws.Range("L136").Value = 2175.8571
$ws.Range("M136").Value = -1359
$ws.Range("N136").Value = -7275.8571
